$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update candidate data rows 2-5 with new values
$ws.Range("A2").Value = "dvjjB737"
$ws.Range("B2").Value = 231004175
$ws.Range("C2").Value = "exbgpip37"
$ws.Range("D2").Value = "X5!r&Gq2"
$ws.Range("F2").Value = "LmnlFMkK"
$ws.Range("G2").Value = "KLfZ"

$ws.Range("A3").Value = "KxbjB928"
$ws.Range("B3").Value = 231004174
$ws.Range("C3").Value = "mldfrft48"
$ws.Range("D3").Value = "Ja9%N7x$"
$ws.Range("F3").Value = "PSgJEXkU"
$ws.Range("G3").Value = "viiC"

$ws.Range("A4").Value = "HAlPw446"
$ws.Range("B4").Value = 231004173
$ws.Range("C4").Value = "iyiiors16"
$ws.Range("D4").Value = "u&2H%pS4"
$ws.Range("F4").Value = "JSyHNNWG"
$ws.Range("G4").Value = "zQlr"

$ws.Range("A5").Value = "vaTQN951"
$ws.Range("B5").Value = 231004172
$ws.Range("C5").Value = "wquueqb39"
$ws.Range("D5").Value = "Yr7&Q8x%"
$ws.Range("F5").Value = "kkGAAPzC"
$ws.Range("G5").Value = "YfsI"

# Delete row 6 entirely (was the 5th candidate row)
$ws.Rows("6:6").Delete()

# Update the selection / active range to match the new data extent
$ws.Range("A1:H5").Select()

